# Apply the updated Leve-profit figures produced by the scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4546328
$ws.Range("J17").Value = 5000910.5
$ws.Range("L17").Value = 15002731.5
$ws.Range("N17").Value = -15003067.5
$ws.Range("H41").Value = 492.5625
$ws.Range("I41").Value = 609
$ws.Range("J41").Value = 298.5
$ws.Range("K41").Value = 609
$ws.Range("L41").Value = 298.5
$ws.Range("M41").Value = -169
$ws.Range("N41").Value = -1178.5
$ws.Range("H53").Value = 351
$ws.Range("I53").Value = 376.25
$ws.Range("J53").Value = 250
$ws.Range("K53").Value = 376.25
$ws.Range("L53").Value = 250
$ws.Range("M53").Value = 260.75
$ws.Range("N53").Value = -1524
$ws.Range("H62").Value = 3107.4614
$ws.Range("I62").Value = 3133.2222
$ws.Range("J62").Value = 3049.5
$ws.Range("K62").Value = 3133.2222
$ws.Range("L62").Value = 3049.5
$ws.Range("M62").Value = -2509.2222
$ws.Range("N62").Value = -4297.5
$ws.Range("H65").Value = 3107.4614
$ws.Range("I65").Value = 3133.2222
$ws.Range("J65").Value = 3049.5
$ws.Range("K65").Value = 15666.111
$ws.Range("L65").Value = 15247.5
$ws.Range("M65").Value = -12546.111
$ws.Range("N65").Value = -21487.5
$ws.Range("H76").Value = 3481.6
$ws.Range("I76").Value = 2962.3057
$ws.Range("J76").Value = 5558.778
$ws.Range("K76").Value = 2962.3057
$ws.Range("L76").Value = 5558.778
$ws.Range("M76").Value = -2647.3057
$ws.Range("N76").Value = -6188.778
$ws.Range("H79").Value = 3481.6
$ws.Range("I79").Value = 2962.3057
$ws.Range("J79").Value = 5558.778
$ws.Range("K79").Value = 2962.3057
$ws.Range("L79").Value = 5558.778
$ws.Range("M79").Value = -1870.3057
$ws.Range("N79").Value = -7742.778
$ws.Range("H86").Value = 23787.889
$ws.Range("I86").Value = 1339.2727
$ws.Range("K86").Value = 1339.2727
$ws.Range("M86").Value = -216.2727
$ws.Range("H89").Value = 23787.889
$ws.Range("I89").Value = 1339.2727
$ws.Range("K89").Value = 6696.363499999999
$ws.Range("M89").Value = -1080.363499999999
$ws.Range("H106").Value = 1024.875
$ws.Range("I106").Value = 1024.875
$ws.Range("K106").Value = 1024.875
$ws.Range("M106").Value = -393.875
$ws.Range("H112").Value = 1321.6522
$ws.Range("J112").Value = 1353.0526
$ws.Range("L112").Value = 4059.1578
$ws.Range("N112").Value = -6275.1578
$ws.Range("H116").Value = 2796.182
$ws.Range("J116").Value = 2651
$ws.Range("L116").Value = 2651
$ws.Range("N116").Value = -9535
$ws.Range("H138").Value = 29279.52
$ws.Range("I138").Value = 1270.027
$ws.Range("J138").Value = 56551.92
$ws.Range("K138").Value = 3810.081
$ws.Range("L138").Value = 169655.76
$ws.Range("M138").Value = 1329.919
$ws.Range("N138").Value = -179935.76

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2641.1875
$ws.Range("I105").Value = 2641.1875
$ws.Range("K105").Value = 2641.1875
$ws.Range("M105").Value = -894.1875

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2540.423
$ws.Range("I58").Value = 2146.3333
$ws.Range("J58").Value = 2878.2144
$ws.Range("K58").Value = 2146.3333
$ws.Range("L58").Value = 2878.2144
$ws.Range("M58").Value = -1943.3333
$ws.Range("N58").Value = -3284.2144
$ws.Range("H99").Value = 23340.4
$ws.Range("I99").Value = 2234
$ws.Range("J99").Value = 55000
$ws.Range("K99").Value = 2234
$ws.Range("L99").Value = 55000
$ws.Range("M99").Value = -736
$ws.Range("N99").Value = -57996
$ws.Range("H105").Value = 2001.375
$ws.Range("I105").Value = 2374.6667
$ws.Range("K105").Value = 2374.6667
$ws.Range("M105").Value = -627.6667000000002
$ws.Range("H126").Value = 23340.4
$ws.Range("I126").Value = 2234
$ws.Range("J126").Value = 55000
$ws.Range("K126").Value = 6702
$ws.Range("L126").Value = 165000
$ws.Range("M126").Value = -4232
$ws.Range("N126").Value = -169940
$ws.Range("H136").Value = 2540.423
$ws.Range("I136").Value = 2146.3333
$ws.Range("J136").Value = 2878.2144
$ws.Range("K136").Value = 6438.999899999999
$ws.Range("L136").Value = 8634.643199999999
$ws.Range("M136").Value = -3888.999899999999
$ws.Range("N136").Value = -13734.6432

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 713
$ws.Range("J34").Value = 1043.1428
$ws.Range("L34").Value = 3129.4284
$ws.Range("N34").Value = -3297.4284
$ws.Range("H40").Value = 445.96
$ws.Range("I40").Value = 125
$ws.Range("J40").Value = 473.86957
$ws.Range("K40").Value = 500
$ws.Range("L40").Value = 1895.47828
$ws.Range("M40").Value = -431
$ws.Range("N40").Value = -2033.47828
$ws.Range("H69").Value = 2128
$ws.Range("I69").Value = 2170.6667
$ws.Range("J69").Value = 2000
$ws.Range("K69").Value = 6512.000100000001
$ws.Range("L69").Value = 6000
$ws.Range("M69").Value = -5701.000100000001
$ws.Range("N69").Value = -7622
$ws.Range("H72").Value = 2128
$ws.Range("I72").Value = 2170.6667
$ws.Range("J72").Value = 2000
$ws.Range("K72").Value = 19536.0003
$ws.Range("L72").Value = 18000
$ws.Range("M72").Value = -15480.0003
$ws.Range("N72").Value = -26112
$ws.Range("H80").Value = 5000.1
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 5000.1
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 15000.3
$ws.Range("N80").Value = -16872.3
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 5000.1
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 5000.1
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 45000.9
$ws.Range("N83").Value = -54360.9
$ws.Range("M83").ClearContents()
$ws.Range("H92").Value = 1501
$ws.Range("I92").Value = 2002
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 6006
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = -4758
$ws.Range("N92").Value = -5496
$ws.Range("H97").Value = 733.3333
$ws.Range("I97").Value = 300
$ws.Range("J97").Value = 820
$ws.Range("K97").Value = 900
$ws.Range("L97").Value = 2460
$ws.Range("M97").Value = -404
$ws.Range("N97").Value = -3452
$ws.Range("H113").Value = 677.94116
$ws.Range("I113").Value = 552.3
$ws.Range("J113").Value = 857.4286
$ws.Range("K113").Value = 1656.9
$ws.Range("L113").Value = 2572.2858
$ws.Range("M113").Value = 513.1000000000001
$ws.Range("N113").Value = -6912.2858

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 273
$ws.Range("I107").Value = 249.21428
$ws.Range("J107").Value = 356.25
$ws.Range("K107").Value = 747.64284
$ws.Range("L107").Value = 1068.75
$ws.Range("M107").Value = 1172.35716
$ws.Range("N107").Value = -4908.75

Write-Host "Edit complete"
